{"js": "// Replace the whole body with the target content:\n// - paragraph 1: original sentence (now ending with a period), split into\n//   several runs (mirrors the proofing-mark split from the authored edit),\n//   font size bumped to 22pt (sz/szCs = 44 half-points).\n// - paragraph 2: a new blank paragraph (22pt).\n// - paragraph 3: \"X + y = 10\" (22pt).\n// - paragraph 4: \"Y=mx + c (it\\u2019s a equation of linear line)\" (22pt).\nconst flatOpcXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\">This is a git repo session. This file is now being added to the </w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>git</w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> repo</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>.</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n  </w:pPr>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>X</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>+</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>y</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>=</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>10</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>Y=mx</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>+</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\">c (it\u2019s </w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:proofErr w:type=\"gramStart\"/>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>a</w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:proofErr w:type=\"gramEnd\"/>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> equation of linear line)</w:t>\n  </w:r>\n</w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst body = context.document.body;\nbody.insertOoxml(flatOpcXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the whole document body content with the target content:\n# - paragraph 1: original sentence (now ending with a period), split into\n#   several runs (mirroring the proofing-mark split from the authored edit),\n#   font size bumped to 22pt (sz/szCs = 44 half-points).\n# - paragraph 2: a new blank paragraph (22pt).\n# - paragraph 3: \"X + y = 10\" (22pt).\n# - paragraph 4: \"Y=mx + c (it\u2019s a equation of linear line)\" (22pt, curly apostrophe).\n\n$d = $word.ActiveDocument\n\n$wordXml = @'\n<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n  <w:pPr>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\">This is a git repo session. This file is now being added to the </w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>git</w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> repo</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>.</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n  </w:pPr>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>X</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>+</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>y</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>=</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>10</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>Y=mx</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>+</w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\">c (it\u2019s </w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:proofErr w:type=\"gramStart\"/>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t>a</w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:proofErr w:type=\"gramEnd\"/>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"44\"/>\n      <w:szCs w:val=\"44\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> equation of linear line)</w:t>\n  </w:r>\n</w:p>\n'@\n\n\n# $d.Content addresses the whole main-story range; InsertXML() replaces\n# that range's contents with the parsed WordprocessingML (mirrors\n# Office.js's body.insertOoxml(xml, \"Replace\")).\n$r = $d.Content\n$r.InsertXML($wordXml)\n"}
